$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.4
$ws.Range("L3").Value = 6.5
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("U3").Value = 4.3
$ws.Range("V3").Value = 1.22
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 17
$ws.Range("AJ3").Value = 7

# Row 4
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 2.2
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2.05
$ws.Range("U4").Value = 3.75
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 1.18
$ws.Range("Y4").Value = 1.53
$ws.Range("Z4").Value = 2.38
$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AD4").Value = 6
$ws.Range("AF4").Value = 11
$ws.Range("AI4").Value = 7
$ws.Range("AJ4").Value = 7.5
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 101
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 81
$ws.Range("AS4").Value = 67

# Row 5
$ws.Range("S5").Value = 2.3
$ws.Range("T5").Value = 1.6

# Row 6
$ws.Range("G6").Value = 2.63
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 3.75
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("U6").Value = 4
$ws.Range("V6").Value = 1.25
$ws.Range("AA6").Value = 2.05
$ws.Range("AB6").Value = 1.7
$ws.Range("AC6").Value = 6.5
$ws.Range("AD6").Value = 11
$ws.Range("AF6").Value = 26
$ws.Range("AM6").Value = 501
$ws.Range("AN6").Value = 7.5
$ws.Range("AO6").Value = 13
$ws.Range("AP6").Value = 12
$ws.Range("AR6").Value = 29

# Row 8
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 2.9
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("U8").Value = 4.4
$ws.Range("V8").Value = 1.2
$ws.Range("AC8").Value = 7.5
